# Update the "Price" (D) and "Volume(1h)" (E) columns for the cryptos list.
# Values that look like plain numbers (e.g. "555.97") are written with a
# leading apostrophe so Excel stores them as text (matching the source
# workbook's inlineStr cells) instead of silently converting them to
# numeric values and dropping formatting like trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.729.20'
$ws.Range("E2").Value = '  -5.76%  '
$ws.Range("D3").Value = '3.280.96'
$ws.Range("E3").Value = '  -6.27%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '''555.97'
$ws.Range("E5").Value = '  -3.80%  '
$ws.Range("D6").Value = '''183.56'
$ws.Range("E6").Value = '  -4.64%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").Value = '''0.587'
$ws.Range("E8").Value = '  -3.98%  '
$ws.Range("D9").Value = '3.273.65'
$ws.Range("E9").Value = '  -6.17%  '
$ws.Range("D10").Value = '''0.184'
$ws.Range("E10").Value = '  -9.92%  '
$ws.Range("E11").Value = '  -6.20%  '
$ws.Range("D12").Value = '''47.27'
$ws.Range("E12").Value = '  -8.25%  '
$ws.Range("E13").Value = '  -7.18%  '
$ws.Range("D14").Value = '''641.73'
$ws.Range("E14").Value = '  -0.48%  '
$ws.Range("D15").Value = '''8.64'
$ws.Range("E15").Value = '  -5.63%  '
$ws.Range("D16").Value = '3.798.56'
$ws.Range("E16").Value = '  -6.42%  '
$ws.Range("D17").Value = '''18.05'
$ws.Range("E17").Value = '  -1.76%  '
$ws.Range("D18").Value = '65.774.42'
$ws.Range("E18").Value = '  -5.70%  '
$ws.Range("E19").Value = '  -3.28%  '
$ws.Range("D20").Value = '3.274.64'
$ws.Range("E20").Value = '  -6.43%  '
$ws.Range("E21").Value = '  -8.78%  '
$ws.Range("D22").Value = '''0.902'
$ws.Range("E22").Value = '  -5.13%  '
$ws.Range("D23").Value = '''18.35'
$ws.Range("E23").Value = '  +1.70%  '
$ws.Range("D24").Value = '''107.82'
$ws.Range("E24").Value = '  +8.87%  '
$ws.Range("D25").Value = '''4.91'
$ws.Range("E25").Value = '  -8.18%  '
$ws.Range("D26").Value = '''3.96'
$ws.Range("E26").Value = '  -7.45%  '
$ws.Range("D28").Value = '''9.55'
$ws.Range("E28").Value = '  -5.28%  '
$ws.Range("D29").Value = '''8.66'
$ws.Range("E29").Value = '  -7.88%  '
$ws.Range("D30").Value = '''30.27'
$ws.Range("E30").Value = '  -7.43%  '
$ws.Range("D31").Value = '''3.92'
$ws.Range("E31").Value = '  -7.31%  '
$ws.Range("D32").Value = '''6.29'
$ws.Range("E32").Value = '  -6.66%  '
$ws.Range("D33").Value = '''11.03'
$ws.Range("E33").Value = '  -5.40%  '
$ws.Range("D34").Value = '''0.104'
$ws.Range("E34").Value = '  -4.90%  '
$ws.Range("D35").Value = '3.775.29'
$ws.Range("E35").Value = '  +1.49%  '
$ws.Range("D36").Value = '''57.50'
$ws.Range("E36").Value = '  -6.56%  '
$ws.Range("E37").Value = '  -0.14%  '
$ws.Range("D38").Value = '''519.88'
$ws.Range("E38").Value = '  -8.23%  '
$ws.Range("D39").Value = '''3.40'
$ws.Range("E39").Value = '  -5.61%  '
$ws.Range("D40").Value = '0.0₃0733'
$ws.Range("E40").Value = '  -7.24%  '
$ws.Range("D41").Value = '''0.130'
$ws.Range("E41").Value = '  -2.21%  '
$ws.Range("D42").Value = '''2.71'
$ws.Range("E42").Value = '  -6.48%  '
$ws.Range("D43").Value = '''3.36'
$ws.Range("E43").Value = '  -18.16%  '
$ws.Range("D44").Value = '''32.93'
$ws.Range("E44").Value = '  -4.15%  '
$ws.Range("E45").Value = '  -10.28%  '
$ws.Range("E46").Value = '  -6.79%  '
$ws.Range("D47").Value = '''3.21'
$ws.Range("E47").Value = '  -5.01%  '
$ws.Range("E48").Value = '  -4.41%  '
$ws.Range("D49").Value = '''2.61'
$ws.Range("E49").Value = '  -8.42%  '
$ws.Range("D50").Value = '''0.998'
$ws.Range("E50").Value = '  +0.04%  '
$ws.Range("E51").Value = '  +2.08%  '
